# Apply the "Add files via upload" revision to the budget sheet:
#  - append a trailing "." to a handful of unit-price strings in column D
#    (rows 7-12) so they match the "XXX Ft." style used elsewhere
#  - add a new row (13) for "Rack szekrény" (rack cabinet), with its link,
#    quantity, unit price and total price
#  - update the grand total in E15 to reflect the newly added item
#  - refresh the UI selection on the budget sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Költésgvetés")
$ws.Activate()

# --- normalise the "Ft" unit-price strings to end with a period ---------
$ws.Range("D7").Value  = "1.688 Ft."
$ws.Range("D8").Value  = "2.090 Ft."
$ws.Range("D9").Value  = "34.990 Ft."
$ws.Range("D10").Value = "78.099Ft."
$ws.Range("D11").Value = "81.999Ft."
$ws.Range("D12").Value = "46.899 Ft."

# --- new row 13: "Rack szekrény" ----------------------------------------
# Copy formatting from the existing rows so the new row matches the
# look of the rest of the table, then fill in the values.
$ws.Range("A4").Copy()
$ws.Range("A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C4").Copy()
$ws.Range("C13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("D5").Copy()
$ws.Range("D13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("E4").Copy()
$ws.Range("E13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A13").Value = "Rack szekrény"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "98.290 Ft."
$ws.Range("E13").Value = "98.290 Ft."

$ws.Hyperlinks.Add($ws.Range("B13"), "https://ipon.hu/shop/termek/lanberg-wf01-6622-10b/1373766?aku=9a6a3cf85308258c0a45ed35cf45651e")

# --- grand total now includes the rack cabinet ---------------------------
$ws.Range("E15").Value = "10.403.856 Ft."

# --- restore the selection that was active when the file was saved -------
$ws.Range("F18").Select()
